# "fix the test by adding back working.xlsx"
# Trims the sample/test data on the first three sheets back down to a single
# data row each (undoing rows that had accumulated), re-points a couple of
# cells at the "p1942234997" tag value, and moves the active sheet/selection
# from "Cost Center Actuals" to "Cost Center Budgets".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Cost Center Actuals
$ws2 = $wb.Worksheets.Item(2)   # Overall Budgets
$ws3 = $wb.Worksheets.Item(3)   # Cost Center Budgets
$ws5 = $wb.Worksheets.Item(5)   # Tag Bugets

# Grab the "p1942234997" tag string and the Arial/10pt number style (the one
# already used on old row 6) up front, before row 6 gets deleted below.
$tagValue = $ws1.Range("A6").Value2
$ws1.Range("B6").Copy()

# ---------------------------------------------------------------------------
# Sheet "Cost Center Budgets": keep only row 2, re-pointed at the tag value
# with the actual column bumped to 12,000; drop old rows 3-4.
# ---------------------------------------------------------------------------
$ws3.Range("C2").PasteSpecial(-4122)      # xlPasteFormats
$ws3.Range("B2").Value = $tagValue
$ws3.Range("D2").Value = 12000

$ws3.Rows.Item(4).Delete()
$ws3.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# Sheet "Cost Center Actuals": keep only row 2 (re-pointed at the tag value,
# matching the number formatting already used on the old row 6), drop the
# old rows 3-7, and leave row 3 as a blank spacer row.
# ---------------------------------------------------------------------------
$ws1.Range("B2").PasteSpecial(-4122)      # xlPasteFormats
$ws1.Range("A2").Value = $tagValue

$ws1.Rows.Item(7).Delete()
$ws1.Rows.Item(6).Delete()
$ws1.Rows.Item(5).Delete()
$ws1.Rows.Item(4).Delete()

$ws1.Range("A3:B3").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "Overall Budgets": keep only row 2, re-pointed at the tag value with
# the budget bumped to 10,000,000; drop old rows 3-4.
# ---------------------------------------------------------------------------
$ws2.Range("B2").Value = $tagValue
$ws2.Range("C2").Value = 10000000

$ws2.Rows.Item(4).Delete()
$ws2.Rows.Item(3).Delete()

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# View state: active sheet moves from "Cost Center Actuals" to
# "Cost Center Budgets", with fresh selections on every retouched sheet.
# ---------------------------------------------------------------------------
$ws1.Range("A2").Select()
$ws2.Range("B8").Select()
$ws5.Range("A14:XFD14").Select()

$ws3.Activate()
$ws3.Range("E8").Select()
